$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.566.45'
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').Value = '2.243.48'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +1.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('E7').Value = '  +0.81%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0806'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.23'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('D14').Value = '2.288.32'
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.834'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '44.187.62'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '0.0₃0957'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '65.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.97'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('E26').Value = '  +4.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.79'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0797'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.92%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.79'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.29'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.42%  '
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('D43').Value = '1.759.38'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.192'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '80.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '99.15'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '69.96'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  +4.31%  '
